$wb = $excel.ActiveWorkbook

# Switch calculation to automatic (removes calcMode="manual" from calcPr)
$excel.Calculation = -4105

# Internationalize the patient names on the PATIENTS sheet (Name column, B2:B9)
$ws = $wb.Worksheets.Item("PATIENTS")
$ws.Activate()

$ws.Range("B2").Value = "John Doe"
$ws.Range("B3").Value = "Steven Doe"
$ws.Range("B4").Value = "Michael Doe"
$ws.Range("B5").Value = "Jane Doe"
$ws.Range("B6").Value = "Sue Doe"
$ws.Range("B7").Value = "Ellen Doe"
$ws.Range("B8").Value = "Anne Smith"
$ws.Range("B9").Value = "John Smith"

# Leave the selection on B1 of the now-active PATIENTS sheet, matching the
# workbook's new view state (tabSelected moves from DATA to PATIENTS).
$ws.Range("B1").Select()
